$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B30").Value = "求职困难"
$ws.Range("B29").Value = "求职成功途径"

$null = $ws.Range("B29").Select()
